$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 846.6667
$ws.Range("I18").Value = 888
$ws.Range("J18").Value = 826
$ws.Range("K18").Value = 888
$ws.Range("L18").Value = 826
$ws.Range("M18").Value = -604
$ws.Range("N18").Value = -1394

$ws.Range("H129").Value = 741694.5600000001
$ws.Range("J129").Value = 1002165.8
$ws.Range("L129").Value = 3006497.4
$ws.Range("N129").Value = -3016497.4

$ws.Range("H132").Value = 2382833.5
$ws.Range("I132").Value = 1836.3055
$ws.Range("J132").Value = 16668816
$ws.Range("K132").Value = 5508.916499999999
$ws.Range("L132").Value = 50006448
$ws.Range("M132").Value = -2978.916499999999
$ws.Range("N132").Value = -50011508

$ws.Range("H135").Value = 1498.85
$ws.Range("I135").Value = 1914.3077
$ws.Range("J135").Value = 727.2857
$ws.Range("K135").Value = 17228.7693
$ws.Range("L135").Value = 6545.571300000001
$ws.Range("M135").Value = -14693.7693
$ws.Range("N135").Value = -11615.5713

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 13333.333
$ws.Range("J23").Value = 13333.333
$ws.Range("L23").Value = 13333.333
$ws.Range("N23").Value = -13851.333

$ws.Range("H32").Value = 1680.23
$ws.Range("I32").Value = 1688.2323
$ws.Range("J32").Value = 888
$ws.Range("K32").Value = 1688.2323
$ws.Range("L32").Value = 888
$ws.Range("M32").Value = -1401.2323
$ws.Range("N32").Value = -1462

$ws.Range("H45").Value = 1183.7142
$ws.Range("I45").Value = 1122.4
$ws.Range("J45").Value = 1217.7778
$ws.Range("K45").Value = 1122.4
$ws.Range("L45").Value = 1217.7778
$ws.Range("M45").Value = -745.4000000000001
$ws.Range("N45").Value = -1971.7778

$ws.Range("H74").Value = 1079.9524
$ws.Range("I74").Value = 1079
$ws.Range("J74").Value = 1085.6666
$ws.Range("K74").Value = 1079
$ws.Range("L74").Value = 1085.6666
$ws.Range("M74").Value = -205
$ws.Range("N74").Value = -2833.6666

$ws.Range("H77").Value = 1079.9524
$ws.Range("I77").Value = 1079
$ws.Range("J77").Value = 1085.6666
$ws.Range("K77").Value = 5395
$ws.Range("L77").Value = 5428.333000000001
$ws.Range("M77").Value = -1027
$ws.Range("N77").Value = -14164.333

$ws.Range("H132").Value = 1086.0492
$ws.Range("I132").Value = 935.36206
$ws.Range("J132").Value = 3999.3333
$ws.Range("K132").Value = 2806.08618
$ws.Range("L132").Value = 11997.9999
$ws.Range("M132").Value = -276.0861800000002
$ws.Range("N132").Value = -17057.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H102").Value = 14500
$ws.Range("I102").Value = 14500
$ws.Range("K102").Value = 14500
$ws.Range("M102").Value = -11255

$ws.Range("H124").Value = 15000
$ws.Range("J124").Value = 15000
$ws.Range("L124").Value = 15000
$ws.Range("N124").Value = -24820

$ws.Range("H134").Value = 27590.46
$ws.Range("I134").Value = 1756.7059
$ws.Range("J134").Value = 203260
$ws.Range("K134").Value = 5270.1177
$ws.Range("L134").Value = 609780
$ws.Range("M134").Value = -2735.1177
$ws.Range("N134").Value = -614850

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 25609.094
$ws.Range("I31").Value = 2263.0667
$ws.Range("J31").Value = 79484.53999999999
$ws.Range("K31").Value = 2263.0667
$ws.Range("L31").Value = 79484.53999999999
$ws.Range("M31").Value = -1968.0667
$ws.Range("N31").Value = -80074.53999999999

$ws.Range("H34").Value = 25609.094
$ws.Range("I34").Value = 2263.0667
$ws.Range("J34").Value = 79484.53999999999
$ws.Range("K34").Value = 2263.0667
$ws.Range("L34").Value = 79484.53999999999
$ws.Range("M34").Value = -2061.0667
$ws.Range("N34").Value = -79888.53999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 16677023
$ws.Range("I68").Value = 20833780
$ws.Range("J68").Value = 50000
$ws.Range("K68").Value = 62501340
$ws.Range("L68").Value = 150000
$ws.Range("M68").Value = -62500529
$ws.Range("N68").Value = -151622

$ws.Range("H71").Value = 16677023
$ws.Range("I71").Value = 20833780
$ws.Range("J71").Value = 50000
$ws.Range("K71").Value = 187504020
$ws.Range("L71").Value = 450000
$ws.Range("M71").Value = -187499964
$ws.Range("N71").Value = -458112

$ws.Range("H113").Value = 714.35297
$ws.Range("I113").Value = 767.5714
$ws.Range("J113").Value = 677.1
$ws.Range("K113").Value = 2302.7142
$ws.Range("L113").Value = 2031.3
$ws.Range("M113").Value = -132.7142000000003
$ws.Range("N113").Value = -6371.3

$ws.Range("H117").Value = 3109.4167
$ws.Range("I117").Value = 842.25
$ws.Range("J117").Value = 4243
$ws.Range("K117").Value = 2526.75
$ws.Range("L117").Value = 12729
$ws.Range("M117").Value = 915.25
$ws.Range("N117").Value = -19613

$ws.Range("H121").Value = 898.1724
$ws.Range("J121").Value = 926.6786
$ws.Range("L121").Value = 2780.0358
$ws.Range("N121").Value = -5400.0358

$ws.Range("H132").Value = 1090.8846
$ws.Range("I132").Value = 885.7143
$ws.Range("J132").Value = 1330.25
$ws.Range("K132").Value = 7971.428699999999
$ws.Range("L132").Value = 11972.25
$ws.Range("M132").Value = -5441.428699999999
$ws.Range("N132").Value = -17032.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 7300.15
$ws.Range("J136").Value = 7300.15
$ws.Range("L136").Value = 21900.45
$ws.Range("N136").Value = -27000.45

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 723058.9399999999
$ws.Range("I40").Value = 919420.4399999999
$ws.Range("J40").Value = 3066.6667
$ws.Range("K40").Value = 919420.4399999999
$ws.Range("L40").Value = 3066.6667
$ws.Range("M40").Value = -919284.4399999999
$ws.Range("N40").Value = -3338.6667

$ws.Range("H122").Value = 8252.941000000001
$ws.Range("I122").Value = 9013.333000000001
$ws.Range("J122").Value = 2550
$ws.Range("K122").Value = 27039.999
$ws.Range("L122").Value = 7650
$ws.Range("M122").Value = -24589.999
$ws.Range("N122").Value = -12550

$ws.Range("H132").Value = 1858.5933
$ws.Range("I132").Value = 1749.9762
$ws.Range("J132").Value = 2126.9412
$ws.Range("K132").Value = 5249.9286
$ws.Range("L132").Value = 6380.823600000001
$ws.Range("M132").Value = -2719.9286
$ws.Range("N132").Value = -11440.8236

$ws.Range("H136").Value = 3752.9
$ws.Range("I136").Value = 2066.875
$ws.Range("J136").Value = 6750.278
$ws.Range("K136").Value = 6200.625
$ws.Range("L136").Value = 20250.834
$ws.Range("M136").Value = -3650.625
$ws.Range("N136").Value = -25350.834

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 66668336
$ws.Range("I122").Value = 71430184
$ws.Range("K122").Value = 214290552
$ws.Range("M122").Value = -214288102

$ws.Range("H126").Value = 794.05884
$ws.Range("I126").Value = 731.2
$ws.Range("J126").Value = 968.6667
$ws.Range("K126").Value = 2193.6
$ws.Range("L126").Value = 2906.0001
$ws.Range("M126").Value = 276.3999999999996
$ws.Range("N126").Value = -7846.0001

$ws.Range("H132").Value = 807.5472
$ws.Range("I132").Value = 750.425
$ws.Range("J132").Value = 983.3077
$ws.Range("K132").Value = 2251.275
$ws.Range("L132").Value = 2949.9231
$ws.Range("M132").Value = 278.7250000000004
$ws.Range("N132").Value = -8009.9231
